# Update "想去人数" (interested-count) values on the 展览 and 全部类型 sheets
# to reflect freshly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value, for each sheet (rows differ slightly between sheets
# because 全部类型 aggregates more rows than 展览 before these entries).
$exhibitUpdates = @{
    2  = 1895
    6  = 2696
    10 = 1567
    17 = 7
    22 = 212
    24 = 1734
    27 = 68
}

$allTypeUpdates = @{
    2  = 1895
    7  = 2696
    11 = 1567
    18 = 7
    23 = 212
    25 = 1734
    28 = 68
}

foreach ($row in $exhibitUpdates.Keys) {
    $sheetExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

foreach ($row in $allTypeUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allTypeUpdates[$row]
}
